$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A7 and A8
$ws.Range("A7").Value = "CLIMATE-PROBLEMS"
$ws.Range("A8").Value = "CLIMATE-IMPACТS"

# Row 9 becomes the comment text that used to be at row 11
$ws.Range("A9").Value = "CLIMATE-ORGANIZATIONS is incorrect it should be none because Beijing is a location hence no tag"

# Row 10 is unchanged (CLIMATE-ORGSANISMS)

# Rows 11-15 shift up from what used to be rows 12-16
$ws.Range("A11").Value = "CLIMATE-PROPERTY"
$ws.Range("A12").Value = "CLIMATE-RESEARCH"
$ws.Range("A13").Value = "LOCATION"
$ws.Range("A14").Value = "no label"
$ws.Range("A15").Value = "not a predefined category"

# Delete the now-obsolete last row (row 16) entirely, shifting rows up
$ws.Rows.Item(16).Delete()
